{"js": "// 1) Remove the \"Date du rapport : ...\" paragraph together with the blank\n//    (single-space) paragraph that immediately follows it.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = paragraphs.items;\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].text.indexOf(\"Date du rapport\") !== -1) {\n    // Delete the following blank paragraph first (if present) so indices\n    // for the date paragraph stay valid, then delete the date paragraph.\n    if (i + 1 < items.length) {\n      items[i + 1].delete();\n    }\n    items[i].delete();\n    break;\n  }\n}\nawait context.sync();\n\n// 2) Update wording in the final paragraph: \"exposition nationale\" ->\n//    \"exposition \u00e0 l'\u00e9chelle nationale\".\nconst results = body.search(\"exposition nationale\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < results.items.length; i++) {\n  results.items[i].insertText(\"exposition \u00e0 l\u2019\u00e9chelle nationale\", \"Replace\");\n}\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) Remove the \"Date du rapport : ...\" paragraph together with the blank\n#    (single-space) paragraph that immediately follows it.\n$targetIndex = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    if ($d.Paragraphs.Item($i).Range.Text -like \"*Date du rapport*\") {\n        $targetIndex = $i\n        break\n    }\n}\n\nif ($targetIndex -gt 0) {\n    $nextIndex = $targetIndex + 1\n    if ($nextIndex -le $d.Paragraphs.Count) {\n        $nextText = $d.Paragraphs.Item($nextIndex).Range.Text.Trim()\n        if ($nextText.Length -eq 0) {\n            $d.Paragraphs.Item($nextIndex).Range.Delete()\n        }\n    }\n    $d.Paragraphs.Item($targetIndex).Range.Delete()\n}\n\n# 2) Update wording in the final paragraph: \"exposition nationale\" ->\n#    \"exposition \u00e0 l'\u00e9chelle nationale\".\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"exposition nationale\"\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"exposition \u00e0 l\u2019\u00e9chelle nationale\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n"}
